# Auto push at 09:14:46
# Duplicate the "P1-4B" worksheet to create a new "P1-4B_2" worksheet at the
# end of the workbook, then rewrite its transaction figures with a new
# problem's data set, and finally restore the selection/active-sheet state.

$wb = $excel.ActiveWorkbook

# --- 1. Duplicate P1-4B -> P1-4B_2, placed after the last sheet ---------
$src = $wb.Worksheets.Item("P1-4B")
$src.Copy($null, $wb.Worksheets.Item($wb.Worksheets.Count))
$new = $wb.Worksheets.Item($wb.Worksheets.Count)
$new.Name = "P1-4B_2"

# --- 2. Overwrite the transaction entries with the new figures ----------
$new.Range("B4").Value = 300000
$new.Range("H4").Formula = "=B4"

$new.Range("B5").Value = -60000
$new.Range("E5").Formula = "=-B5"
$new.Range("J5").ClearContents()

$new.Range("B6").Value = -5000
$new.Range("D6").ClearContents()
$new.Range("G6").ClearContents()
$new.Range("J6").Value = -5000

$new.Range("B7").ClearContents()
$new.Range("D7").Value = 500
$new.Range("G7").Value = 500
$new.Range("J7").ClearContents()

$new.Range("B8").Value = -50000
$new.Range("E8").Value = 50000
$new.Range("I8").ClearContents()

$new.Range("B9").Value = 15000
$new.Range("I9").Value = 15000
$new.Range("K9").ClearContents()

$new.Range("B10").Value = -1000
$new.Range("C10").ClearContents()
$new.Range("I10").ClearContents()
$new.Range("J10").Value = -1000

$new.Range("B11").Value = -1200
$new.Range("J11").Value = -1200

$new.Range("B12").ClearContents()
$new.Range("C12").Value = 8000
$new.Range("G12").ClearContents()
$new.Range("I12").Value = 8000

$new.Range("B13").Value = -200
$new.Range("C13").ClearContents()
$new.Range("G13").Value = -200

$new.Range("B14").ClearContents()
$new.Range("G14").ClearContents()

$new.Range("E15").ClearContents()
$new.Range("G15").ClearContents()

$new.Range("B16").ClearContents()
$new.Range("J16").ClearContents()

# --- 3. Selection state on the new sheet ---------------------------------
$new.Range("P16").Select()

# --- 4. Restore per-sheet selections / active sheet ----------------------
$ws1 = $wb.Worksheets.Item("P1-1B")
$ws1.Range("B4").Select()

$ws2 = $wb.Worksheets.Item("P1-2B")
$ws2.Activate()
$ws2.Range("A2").Select()
